# Update income tax notes on the "Joint Tax Assessment for Idris and his wife" slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# --- Table 1 ("Content Placeholder 3" - reliefs / taxable income table) ---
$tbl1 = $s.Shapes.Item(1).Table

# Wife relief: 9,000 -> 4,000
$tbl1.Cell(6,2).Shape.TextFrame.TextRange.Text = "4,000"

# Total Reliefs: (42,600) -> (37,600)
$tbl1.Cell(11,4).Shape.TextFrame.TextRange.Text = "(37,600)"

# Taxable Income: 41,600 -> 46,600
$tbl1.Cell(12,4).Shape.TextFrame.TextRange.Text = "46,600"

# --- Table 2 ("Table 2" - tax rate / payable tax table) ---
$tbl2 = $s.Shapes.Item(2).Table

# "On the next 6,600 x 8%" -> "On the next 11,600 x 8%"
$tbl2.Cell(3,1).Shape.TextFrame.TextRange.Text = "On the next 11,600 x 8%"

# 528 -> 928
$tbl2.Cell(3,2).Shape.TextFrame.TextRange.Text = "928"

# Total Tax: 1,128 -> 1,528
$tbl2.Cell(4,4).Shape.TextFrame.TextRange.Text = "1,528"

# Payable Tax: RM 628 -> RM 1,028
$tbl2.Cell(8,4).Shape.TextFrame.TextRange.Text = "RM 1,028"

# --- TextBox 4: "Balance = 41,600 - 35,000 = 6,600 @ 8% " explanation ---
$tb = $s.Shapes.Item(3).TextFrame.TextRange
$balancePara = $tb.Paragraphs(3,1)
$balanceRun = $balancePara.Runs(1,1)
$balanceRun.Text = "Balance = 46,600 – 35,000 = 11,600 @ 8% "
